$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.767.86"
$ws.Range("E2").Value = "  -3.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.911.20"
$ws.Range("E3").Value = "  -3.97%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.63"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.98"
$ws.Range("E6").Value = "  -6.31%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -2.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.911.63"
$ws.Range("E9").Value = "  -3.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.81"
$ws.Range("E10").Value = "  +5.82%  "
$ws.Range("E11").Value = "  -4.18%  "
$ws.Range("E12").Value = "  -4.19%  "
$ws.Range("E13").Value = "  -3.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.53"
$ws.Range("E14").Value = "  -5.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.395.24"
$ws.Range("E16").Value = "  -3.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.735.74"
$ws.Range("E17").Value = "  -3.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.72"
$ws.Range("E18").Value = "  -5.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.912.01"
$ws.Range("E19").Value = "  -3.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "430.60"
$ws.Range("E20").Value = "  -4.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.61"
$ws.Range("E21").Value = "  -4.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.682"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("E23").Value = "  -4.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.36"
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.79"
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.20"
$ws.Range("E26").Value = "  -4.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.88"
$ws.Range("E27").Value = "  -3.47%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  -3.98%  "
$ws.Range("E31").Value = "  -3.15%  "
$ws.Range("E32").Value = "  -1.70%  "
$ws.Range("E33").Value = "  -3.98%  "
$ws.Range("E34").Value = "  -3.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0865"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  -2.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.66"
$ws.Range("E37").Value = "  -4.30%  "
$ws.Range("E38").Value = "  -1.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("E39").Value = "  -5.45%  "
$ws.Range("E40").Value = "  -2.10%  "
$ws.Range("E41").Value = "  -5.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.66"
$ws.Range("E42").Value = "  -4.78%  "
$ws.Range("E43").Value = "  -3.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.33"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("E45").Value = "  -3.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "373.33"
$ws.Range("E46").Value = "  -5.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.669.86"
$ws.Range("E47").Value = "  -2.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.90"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.30"
$ws.Range("E51").Value = "  -1.83%  "
